$d = $word.ActiveDocument

# Build the paragraph's new content (three runs with a spell-check
# proofErr wrapper around "segração") as a WordProcessingML package
# fragment and drop it into the whole document story. Word resolves a
# package-wrapped InsertXML against the full story range as a content
# replace rather than an insert-before, so this swaps the lone empty
# paragraph's content in place instead of leaving a stray empty
# paragraph behind.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Para o princípio da </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>segração</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> da interface, poderia ter sido utilizado uma interface de pagamento que permitisse a comunicação com qualquer sistema de validação de compras, sem prejuízo da completude das transações. Assim, foi criada uma interface de pagamento para realizar esta comunicação.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml) | Out-Null
